$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (border + bold + centered) from H1 to the new header cells I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I (I0) and J (IF)
$data = @(
    @(6, 6),
    @(7, 7),
    @(6, 7),
    @(7, 7),
    @(7, 7),
    @(8, 9),
    @(6, 7),
    @(3, 4),
    @(7, 7),
    @(7, 7),
    @(9, 9),
    @(6, 7),
    @(3, 5),
    @(12, 12),
    @(5, 6),
    @(8, 8),
    @(6, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
